# Update the division-problem worksheet table: each populated row of the
# table (rows 1, 5, 9, 13, 17 — the others are blank spacer rows) gets its
# five cell values replaced with the new problem set, cell by cell, so
# that run formatting (font / size) on each cell is left untouched.

$d = $word.ActiveDocument
$t = $d.Tables(1)

# Row => ordered list of new cell texts (5 columns).
$updates = @{
    1  = @("59÷9=", "91÷4=", "23÷4=", "59÷3=", "41÷6=")
    5  = @("44÷2=", "56÷3=", "87÷9=", "87÷4=", "60÷8=")
    9  = @("96÷5=", "81÷9=", "85÷3=", "98÷2=", "78÷4=")
    13 = @("61÷5=", "54÷7=", "60÷4=", "36÷3=", "37÷2=")
    17 = @("47÷7=", "69÷5=", "10÷7=", "40÷2=", "33÷3=")
}

foreach ($rowIndex in $updates.Keys) {
    $values = $updates[$rowIndex]
    $row = $t.Rows($rowIndex)
    for ($col = 1; $col -le $values.Length; $col++) {
        $row.Cells($col).Range.Text = $values[$col - 1]
    }
}
